$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.0292345
$ws.Range("H2").Value = 0.058469
$ws.Range("I2").Value = 0.4428765120700495
$ws.Range("J2").Value = 0.346386487911515
$ws.Range("M2").Value = 12.3291175
$ws.Range("N2").Value = 24.658235
$ws.Range("O2").Value = 0.1991607983368005
$ws.Range("P2").Value = 0.1614453197874725
$ws.Range("Q2").Value = 0.36043558555375
$ws.Range("R2").Value = 1.441742342215
$ws.Range("S2").Value = 0.08820363970848873
$ws.Range("T2").Value = 0.05592247731093401
$ws.Range("G3").Value = 0.0292345
$ws.Range("H3").Value = 0.058469
$ws.Range("I3").Value = 0.4428765120700495
$ws.Range("J3").Value = 0.346386487911515
$ws.Range("O3").Value = 0.4636226915653649
$ws.Range("P3").Value = 0.563738303362699
$ws.Range("Q3").Value = 0.8390512475641667
$ws.Range("R3").Value = 5.034307485385001
$ws.Range("S3").Value = 0.2053276005569971
$ws.Range("T3").Value = 0.1952713310030015
$ws.Range("G4").Value = 0.0292345
$ws.Range("H4").Value = 0.058469
$ws.Range("I4").Value = 0.4428765120700495
$ws.Range("J4").Value = 0.346386487911515
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.073231
$ws.Range("N4").Value = 0.219693
$ws.Range("O4").Value = 0.001182951206605196
$ws.Range("P4").Value = 0.001438400057427841
$ws.Range("Q4").Value = 0.0021408716695
$ws.Range("R4").Value = 0.012845230017
$ws.Range("S4").Value = 0.0005239013043303658
$ws.Range("T4").Value = 0.0004982423441041513
$ws.Range("G5").Value = 0.0292345
$ws.Range("H5").Value = 0.058469
$ws.Range("I5").Value = 0.4428765120700495
$ws.Range("J5").Value = 0.346386487911515
$ws.Range("M5").Value = 20.652629
$ws.Range("N5").Value = 41.30525799999999
$ws.Range("O5").Value = 0.3336162608064818
$ws.Range("P5").Value = 0.2704386825218454
$ws.Range("Q5").Value = 0.6037692825004999
$ws.Range("R5").Value = 2.415077130002
$ws.Range("S5").Value = 0.1477508059558266
$ws.Range("T5").Value = 0.09367630543415924
$ws.Range("G6").Value = 0.0292345
$ws.Range("H6").Value = 0.058469
$ws.Range("I6").Value = 0.4428765120700495
$ws.Range("J6").Value = 0.346386487911515
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.05032066666666666
$ws.Range("N6").Value = 0.150962
$ws.Range("O6").Value = 0.0008128646795825703
$ws.Range("P6").Value = 0.0009883963051595711
$ws.Range("Q6").Value = 0.001471099529666666
$ws.Range("R6").Value = 0.008826597177999999
$ws.Range("S6").Value = 0.0003599986740784671
$ws.Range("T6").Value = 0.0003423671248089419
$ws.Range("G7").Value = 0.0292345
$ws.Range("H7").Value = 0.058469
$ws.Range("I7").Value = 0.4428765120700495
$ws.Range("J7").Value = 0.346386487911515
$ws.Range("M7").Value = 0.09932299999999999
$ws.Range("N7").Value = 0.297969
$ws.Range("O7").Value = 0.001604433405165134
$ws.Range("P7").Value = 0.001950897965395876
$ws.Range("Q7").Value = 0.0029036582435
$ws.Range("R7").Value = 0.017421949461
$ws.Range("S7").Value = 0.0007105658703282069
$ws.Range("T7").Value = 0.0006757646945071979
$ws.Range("I8").Value = 0.5571234879299505
$ws.Range("J8").Value = 0.6536135120884849
$ws.Range("M8").Value = 12.3291175
$ws.Range("N8").Value = 24.658235
$ws.Range("O8").Value = 0.1991607983368005
$ws.Range("P8").Value = 0.1614453197874725
$ws.Range("Q8").Value = 0.4534156251799999
$ws.Range("R8").Value = 2.72049375108
$ws.Range("S8").Value = 0.1109571586283118
$ws.Range("T8").Value = 0.1055228424765385
$ws.Range("I9").Value = 0.5571234879299505
$ws.Range("J9").Value = 0.6536135120884849
$ws.Range("O9").Value = 0.4636226915653649
$ws.Range("P9").Value = 0.563738303362699
$ws.Range("S9").Value = 0.2582950910083677
$ws.Range("T9").Value = 0.3684669723596975
$ws.Range("I10").Value = 0.5571234879299505
$ws.Range("J10").Value = 0.6536135120884849
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 0.6666666666666666
$ws.Range("M10").Value = 0.073231
$ws.Range("N10").Value = 0.219693
$ws.Range("O10").Value = 0.001182951206605196
$ws.Range("P10").Value = 0.001438400057427841
$ws.Range("Q10").Value = 0.002693143256
$ws.Range("R10").Value = 0.024238289304
$ws.Range("S10").Value = 0.0006590499022748305
$ws.Range("T10").Value = 0.0009401577133236896
$ws.Range("I11").Value = 0.5571234879299505
$ws.Range("J11").Value = 0.6536135120884849
$ws.Range("M11").Value = 20.652629
$ws.Range("N11").Value = 41.30525799999999
$ws.Range("O11").Value = 0.3336162608064818
$ws.Range("P11").Value = 0.2704386825218454
$ws.Range("Q11").Value = 0.7595210841039999
$ws.Range("R11").Value = 4.557126504623999
$ws.Range("S11").Value = 0.1858654548506551
$ws.Range("T11").Value = 0.1767623770876861
$ws.Range("I12").Value = 0.5571234879299505
$ws.Range("J12").Value = 0.6536135120884849
$ws.Range("K12").Value = 2
$ws.Range("L12").Value = 0.6666666666666666
$ws.Range("M12").Value = 0.05032066666666666
$ws.Range("N12").Value = 0.150962
$ws.Range("O12").Value = 0.0008128646795825703
$ws.Range("P12").Value = 0.0009883963051595711
$ws.Range("Q12").Value = 0.001850592837333333
$ws.Range("R12").Value = 0.016655335536
$ws.Range("S12").Value = 0.0004528660055041032
$ws.Range("T12").Value = 0.0006460291803506292
$ws.Range("I13").Value = 0.5571234879299505
$ws.Range("J13").Value = 0.6536135120884849
$ws.Range("M13").Value = 0.09932299999999999
$ws.Range("N13").Value = 0.297969
$ws.Range("O13").Value = 0.001604433405165134
$ws.Range("P13").Value = 0.001950897965395876
$ws.Range("Q13").Value = 0.003652702647999999
$ws.Range("R13").Value = 0.03287432383199999
$ws.Range("S13").Value = 0.0008938675348369266
$ws.Range("T13").Value = 0.001275133270888678
